$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - reorder the category labels
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "kitchens_1"
$ws.Range("F1").Value = "kitchens_2"

# Data rows 2-7 - updated one-hot selection values
$data = @(
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0),
    @(0,0,0,0,0,1),
    @(0,1,0,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
